# QuestionnaireAnnexes.xlsx edit:
#  - rename annex sheet tabs to add spaces around the dash
#  - fix the Print_Titles defined name to match the renamed sheet
#  - add a new "ANN 5 - Staff Codes" sheet with staff roster data
#  - move the active/selected tab from ANN4 to ANN1

$wb = $excel.ActiveWorkbook

# --- 1. Rename the first four annex tabs -------------------------------
$wb.Worksheets.Item(1).Name = "ANN 1 - Province Codes"
$wb.Worksheets.Item(2).Name = "ANN 2 - District Codes"
$wb.Worksheets.Item(3).Name = "ANN 3 - Birthplace Codes"
$wb.Worksheets.Item(4).Name = "ANN 4 - Possession Value Limits"

# --- 2. Fix the Print_Titles defined name so it points at the renamed sheet
foreach ($n in $wb.Names) {
    if ($n.Name -like "*Print_Titles*") {
        $n.RefersTo = "='ANN 3 - Birthplace Codes'!`$1:`$1"
    }
}

# --- 3. Add the new "ANN 5 - Staff Codes" sheet after ANN 4 ------------
$lastIndex = $wb.Worksheets.Count
$ws5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($lastIndex))
$ws5.Name = "ANN 5 - Staff Codes"

# Column widths (approximate the original character based widths)
$ws5.Columns.Item(1).ColumnWidth = 8.498697916666666
$ws5.Columns.Item(2).ColumnWidth = 19.608072916666668
$ws5.Columns.Item(3).ColumnWidth = 13.830729166666666
$ws5.Columns.Item(4).ColumnWidth = 7.166666666666667
$ws5.Columns.Item(5).ColumnWidth = 5.944010416666667
$ws5.Columns.Item(6).ColumnWidth = 2.2760416666666665

# Header row
$ws5.Range("A1").Value = "Staff code"
$ws5.Range("B1").Value = "Name"
$ws5.Range("C1").Value = "Role (1=interviewer, 2=supervisor)"
$ws5.Range("D1").Value = "Province"
$ws5.Range("E1").Value = "District"
$ws5.Range("F1").Value = "EA"

# Staff roster rows
$ws5.Range("A2").Value = 1
$ws5.Range("B2").Value = "Shemika Rothenberger  "
$ws5.Range("C2").Value = 2
$ws5.Range("D2").Value = 1
$ws5.Range("E2").Value = 1

$ws5.Range("A3").Value = 2
$ws5.Range("B3").Value = "Andrew Benninger  "
$ws5.Range("C3").Value = 1
$ws5.Range("D3").Value = 1
$ws5.Range("E3").Value = 1
$ws5.Range("F3").Value = 1

$ws5.Range("A4").Value = 3
$ws5.Range("B4").Value = "Angelica Swenson  "
$ws5.Range("C4").Value = 1
$ws5.Range("D4").Value = 1
$ws5.Range("E4").Value = 1
$ws5.Range("F4").Value = 2

$ws5.Range("A5").Value = 4
$ws5.Range("B5").Value = "Zelma Hawke  "
$ws5.Range("C5").Value = 1
$ws5.Range("D5").Value = 1
$ws5.Range("E5").Value = 1
$ws5.Range("F5").Value = 3

$ws5.Range("A6").Value = 5
$ws5.Range("B6").Value = "Willis Catron "
$ws5.Range("C6").Value = 1
$ws5.Range("D6").Value = 1
$ws5.Range("E6").Value = 1
$ws5.Range("F6").Value = 4

# Header row formatting: shaded fill, thin box border, wrapped text
$header = $ws5.Range("A1:F1")
$header.RowHeight = 43.2
$header.WrapText = $true
$header.Interior.Color = 14806254
$header.Borders.LineStyle = 1

# Staff code column: "000" number format with thin border
$ws5.Range("A2:A6").NumberFormat = "000"
$ws5.Range("A2:A6").Borders.LineStyle = 1

# Remaining data cells: thin border
$ws5.Range("B2:F6").Borders.LineStyle = 1

$ws5.Range("I20").Select()

# --- 4. Move the active tab from ANN4 back to ANN1 ----------------------
$wb.Worksheets.Item(1).Activate()
